$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 95.92308
$ws.Range("I9").Value = 95.75
$ws.Range("K9").Value = 95.75
$ws.Range("M9").Value = 73.25
$ws.Range("H12").Value = 173.71428
$ws.Range("I12").Value = 173.71428
$ws.Range("K12").Value = 173.71428
$ws.Range("M12").Value = -3.714280000000002
$ws.Range("H17").Value = 974.3421
$ws.Range("J17").Value = 946.8276
$ws.Range("L17").Value = 2840.4828
$ws.Range("N17").Value = -3176.4828
$ws.Range("H74").Value = 8722.682000000001
$ws.Range("I74").Value = 8899.951999999999
$ws.Range("K74").Value = 8899.951999999999
$ws.Range("M74").Value = -7963.951999999999
$ws.Range("H77").Value = 8722.682000000001
$ws.Range("I77").Value = 8899.951999999999
$ws.Range("K77").Value = 44499.75999999999
$ws.Range("M77").Value = -39819.75999999999
$ws.Range("H106").Value = 7586
$ws.Range("I106").Value = 4200.364
$ws.Range("J106").Value = 20000
$ws.Range("K106").Value = 4200.364
$ws.Range("L106").Value = 20000
$ws.Range("M106").Value = -3569.364
$ws.Range("N106").Value = -21262
$ws.Range("H113").Value = 4250
$ws.Range("I113").Value = 4250
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4250
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H116").Value = 28829.8
$ws.Range("I116").Value = 29750
$ws.Range("K116").Value = 29750
$ws.Range("M116").Value = -26308
$ws.Range("H125").Value = 4000
$ws.Range("I125").Value = 2000
$ws.Range("J125").Value = 6000
$ws.Range("K125").Value = 18000
$ws.Range("L125").Value = 54000
$ws.Range("M125").Value = -15540
$ws.Range("N125").Value = -58920
$ws.Range("H132").Value = 13343.0625
$ws.Range("I132").Value = 15170.643
$ws.Range("K132").Value = 45511.929
$ws.Range("M132").Value = -42981.929
$ws.Range("H137").Value = 4120.9287
$ws.Range("I137").Value = 1766.1666
$ws.Range("K137").Value = 5298.4998
$ws.Range("M137").Value = -2748.4998

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 828.619
$ws.Range("I2").Value = 850.17645
$ws.Range("K2").Value = 850.17645
$ws.Range("M2").Value = -737.17645
$ws.Range("H61").Value = 3184.8572
$ws.Range("I61").Value = 3215.6667
$ws.Range("K61").Value = 3215.6667
$ws.Range("M61").Value = -3003.6667
$ws.Range("H88").Value = 2462.7917
$ws.Range("J88").Value = 2462.7917
$ws.Range("L88").Value = 2462.7917
$ws.Range("N88").Value = -3274.7917
$ws.Range("H91").Value = 2462.7917
$ws.Range("J91").Value = 2462.7917
$ws.Range("L91").Value = 2462.7917
$ws.Range("N91").Value = -5270.7917
$ws.Range("H116").Value = 828.619
$ws.Range("I116").Value = 850.17645
$ws.Range("K116").Value = 850.17645
$ws.Range("M116").Value = 1443.82355
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H136").Value = 3184.8572
$ws.Range("I136").Value = 3215.6667
$ws.Range("K136").Value = 9647.000100000001
$ws.Range("M136").Value = -7097.000100000001

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 828.619
$ws.Range("I3").Value = 850.17645
$ws.Range("K3").Value = 850.17645
$ws.Range("M3").Value = -736.17645
$ws.Range("H80").Value = 447.33334
$ws.Range("I80").Value = 503.9091
$ws.Range("J80").Value = 358.42856
$ws.Range("K80").Value = 503.9091
$ws.Range("L80").Value = 358.42856
$ws.Range("M80").Value = 494.0909
$ws.Range("N80").Value = -2354.42856
$ws.Range("H83").Value = 447.33334
$ws.Range("I83").Value = 503.9091
$ws.Range("J83").Value = 358.42856
$ws.Range("K83").Value = 2519.5455
$ws.Range("L83").Value = 1792.1428
$ws.Range("M83").Value = 2472.4545
$ws.Range("N83").Value = -11776.1428
$ws.Range("H99").Value = 7434.875
$ws.Range("I99").Value = 7863.8667
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 7863.8667
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = -6365.8667
$ws.Range("N99").Value = -3996
$ws.Range("H122").Value = 68890
$ws.Range("J122").Value = 68890
$ws.Range("L122").Value = 68890
$ws.Range("N122").Value = -78690
$ws.Range("H134").Value = 4615.6206
$ws.Range("I134").Value = 5050.478
$ws.Range("K134").Value = 15151.434
$ws.Range("M134").Value = -12616.434

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1673.375
$ws.Range("I22").Value = 1070.8182
$ws.Range("J22").Value = 2999
$ws.Range("K22").Value = 1070.8182
$ws.Range("L22").Value = 2999
$ws.Range("M22").Value = -720.8181999999999
$ws.Range("N22").Value = -3699
$ws.Range("H31").Value = 2065
$ws.Range("I31").Value = 2065
$ws.Range("K31").Value = 2065
$ws.Range("M31").Value = -1770
$ws.Range("H34").Value = 2065
$ws.Range("I34").Value = 2065
$ws.Range("K34").Value = 2065
$ws.Range("M34").Value = -1863
$ws.Range("H58").Value = 2911.5625
$ws.Range("J58").Value = 3552.111
$ws.Range("L58").Value = 3552.111
$ws.Range("N58").Value = -3958.111
$ws.Range("H62").Value = 4233.1665
$ws.Range("I62").Value = 3850
$ws.Range("K62").Value = 3850
$ws.Range("M62").Value = -3226
$ws.Range("H65").Value = 4233.1665
$ws.Range("I65").Value = 3850
$ws.Range("K65").Value = 19250
$ws.Range("M65").Value = -16130
$ws.Range("H105").Value = 14036.182
$ws.Range("I105").Value = 15139.8
$ws.Range("K105").Value = 15139.8
$ws.Range("M105").Value = -13392.8
$ws.Range("H132").Value = 22572.25
$ws.Range("I132").Value = 28763
$ws.Range("K132").Value = 86289
$ws.Range("M132").Value = -83759
$ws.Range("H136").Value = 2911.5625
$ws.Range("J136").Value = 3552.111
$ws.Range("L136").Value = 10656.333
$ws.Range("N136").Value = -15756.333

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H37").Value = 42499.445
$ws.Range("J37").Value = 42499.445
$ws.Range("L37").Value = 127498.335
$ws.Range("N37").Value = -127722.335
$ws.Range("H54").Value = 2099
$ws.Range("J54").Value = 2099
$ws.Range("L54").Value = 6297
$ws.Range("N54").Value = -7415
$ws.Range("H68").Value = 3160
$ws.Range("J68").Value = 4000
$ws.Range("L68").Value = 12000
$ws.Range("N68").Value = -13622
$ws.Range("H71").Value = 3160
$ws.Range("J71").Value = 4000
$ws.Range("L71").Value = 36000
$ws.Range("N71").Value = -44112
$ws.Range("H122").Value = 2020252.5
$ws.Range("J122").Value = 4666.3335
$ws.Range("L122").Value = 41997.0015
$ws.Range("N122").Value = -46897.0015
$ws.Range("H128").Value = 316378.3
$ws.Range("I128").Value = 316378.3
$ws.Range("K128").Value = 949134.8999999999
$ws.Range("M128").Value = -944154.8999999999

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 5333.3335
$ws.Range("I3").Value = 3500
$ws.Range("J3").Value = 6250
$ws.Range("K3").Value = 3500
$ws.Range("L3").Value = 6250
$ws.Range("M3").Value = -3384
$ws.Range("N3").Value = -6482
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H11").Value = 27061050
$ws.Range("I11").Value = 8704490
$ws.Range("J11").Value = 50006750
$ws.Range("K11").Value = 8704490
$ws.Range("L11").Value = 50006750
$ws.Range("M11").Value = -8704351
$ws.Range("N11").Value = -50007028
$ws.Range("H14").Value = 7147572
$ws.Range("I14").Value = 10003401
$ws.Range("K14").Value = 10003401
$ws.Range("M14").Value = -10003233
$ws.Range("H126").Value = 2699.3572
$ws.Range("I126").Value = 2531.2222
$ws.Range("J126").Value = 3002
$ws.Range("K126").Value = 7593.6666
$ws.Range("L126").Value = 9006
$ws.Range("M126").Value = -5123.6666
$ws.Range("N126").Value = -13946
$ws.Range("H132").Value = 12535.934
$ws.Range("I132").Value = 13767.704
$ws.Range("J132").Value = 1450
$ws.Range("K132").Value = 41303.112
$ws.Range("L132").Value = 4350
$ws.Range("M132").Value = -38773.112
$ws.Range("N132").Value = -9410

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 9433.695
$ws.Range("I61").Value = 8049
$ws.Range("K61").Value = 8049
$ws.Range("M61").Value = -7847
$ws.Range("H113").Value = 9433.695
$ws.Range("I113").Value = 8049
$ws.Range("K113").Value = 8049
$ws.Range("M113").Value = -5879

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 1336664.4
$ws.Range("I4").Value = 1999996.5
$ws.Range("J4").Value = 10000
$ws.Range("K4").Value = 1999996.5
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = -1999883.5
$ws.Range("N4").Value = -10226
$ws.Range("H126").Value = 3101.3333
$ws.Range("I126").Value = 1004
$ws.Range("K126").Value = 3012
$ws.Range("M126").Value = -542
$ws.Range("H132").Value = 2573.75
$ws.Range("I132").Value = 1938.9286
$ws.Range("K132").Value = 5816.7858
$ws.Range("M132").Value = -3286.7858
